# Bug of histogram fix
# Rewrites the per-child histogram rows (6-10), collapses the summary rows
# (school/cost/time) up into rows 11-13, and drops the now-unused rows 14-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
}

# --- nChildren summary (row 4) ---
Set-TextCell 4 2 "5"

# --- row 6 ---
Set-TextCell 6 2 "2"
$ws.Cells.Item(6,3).Value2 = "Elwanda " + [char]160
$ws.Cells.Item(6,4).Value2 = "Cassy " + [char]160
$ws.Cells.Item(6,5).Value2 = "5,9"
$ws.Cells.Item(6,6).Value2 = "Tamisha(mother): 0550693864"
Set-TextCell 6 8 "21.0"

# --- row 7 ---
Set-TextCell 7 2 "0"
$ws.Cells.Item(7,3).Value2 = "Trudie " + [char]160
$ws.Cells.Item(7,4).Value2 = "Fleta " + [char]160
$ws.Cells.Item(7,5).Value2 = "7,9"
$ws.Cells.Item(7,6).Value2 = "Anneliese(father): 0548973345"
$ws.Cells.Item(7,7).Value2 = "7:03:00"
Set-TextCell 7 8 "18.0"

# --- row 8 ---
Set-TextCell 8 2 "4"
$ws.Cells.Item(8,3).Value2 = "Francisca " + [char]160
$ws.Cells.Item(8,4).Value2 = "Stevie " + [char]160
$ws.Cells.Item(8,5).Value2 = "8,7"
$ws.Cells.Item(8,6).Value2 = "Bernardine(mother): 0561339273"
$ws.Cells.Item(8,7).Value2 = "7:06:00"
Set-TextCell 8 8 "15.0"

# --- row 9 ---
Set-TextCell 9 2 "3"
$ws.Cells.Item(9,3).Value2 = "Alexia " + [char]160
$ws.Cells.Item(9,4).Value2 = "Ramonita " + [char]160
$ws.Cells.Item(9,5).Value2 = "7,6"
$ws.Cells.Item(9,6).Value2 = "Han(father): 0567537032"
$ws.Cells.Item(9,7).Value2 = "7:08:00"
Set-TextCell 9 8 "13.0"

# --- row 10 ---
Set-TextCell 10 2 "1"
$ws.Cells.Item(10,3).Value2 = "Corene " + [char]160
$ws.Cells.Item(10,4).Value2 = "Myra " + [char]160
$ws.Cells.Item(10,5).Value2 = "8,6"
$ws.Cells.Item(10,6).Value2 = "Georgie(mother): 0544823581"
$ws.Cells.Item(10,7).Value2 = "7:10:00"
Set-TextCell 10 8 "11.0"

# --- row 11 becomes the "school" summary row (was row 15) ---
$ws.Cells.Item(11,1).Value2 = "school"
Set-TextCell 11 2 "3"
$ws.Cells.Item(11,3).Value2 = "Ironiah"
$ws.Cells.Item(11,4).Value2 = "mySchool"
$ws.Cells.Item(11,5).Value2 = "0,0"
$ws.Cells.Item(11,6).Value2 = "Shir(secretary): 0523345098"
$ws.Cells.Item(11,7).Value2 = "7:21:00"
$ws.Range("H11").Clear()

# --- row 12 becomes the "cost" summary row (was row 16) ---
$ws.Cells.Item(12,1).Value2 = "cost"
Set-TextCell 12 2 "55.0"
$ws.Range("C12:H12").Clear()

# --- row 13 becomes the "time" summary row (was row 17) ---
$ws.Cells.Item(13,1).Value2 = "time"
Set-TextCell 13 2 "21.0"
$ws.Range("C13:H13").Clear()

# --- rows 14-17 no longer exist ---
$ws.Range("A14:H17").Clear()
